$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 2818
$ws.Cells.Item(2, 5).Value = -361
$ws.Cells.Item(2, 6).Value = -361
$ws.Cells.Item(2, 7).Value = -564
$ws.Cells.Item(2, 8).Value = -564
$ws.Cells.Item(2, 9).Value = -422
$ws.Cells.Item(2, 10).Value = -142
$ws.Cells.Item(2, 11).Value = 8089
$ws.Cells.Item(2, 12).Value = 5382
$ws.Cells.Item(2, 13).Value = 2707
$ws.Cells.Item(2, 14).Value = 2268
$ws.Cells.Item(2, 15).Value = 439
$ws.Cells.Item(2, 16).Value = 655
$ws.Cells.Item(2, 17).Value = 581
$ws.Cells.Item(2, 18).Value = -527
$ws.Cells.Item(2, 19).Value = 9
$ws.Cells.Item(2, 20).Value = 507
$ws.Cells.Item(2, 21).Value = 73
$ws.Cells.Item(2, 22).Value = 3397
$ws.Cells.Item(2, 23).Value = -12.82
$ws.Cells.Item(2, 24).Value = -20.02
$ws.Cells.Item(2, 25).Value = -17.95
$ws.Cells.Item(2, 26).Value = -7.17
$ws.Cells.Item(2, 27).Value = 198.83
$ws.Cells.Item(2, 28).Value = -12.79
$ws.Cells.Item(2, 29).Value = -2736
$ws.Cells.Item(2, 30).Value = -1.56
$ws.Cells.Item(2, 31).Value = 14713
$ws.Cells.Item(2, 32).Value = 0.29
$ws.Cells.Item(2, 33).Value = 0
$ws.Cells.Item(2, 34).Value = 0
$ws.Cells.Item(2, 35).Value = 0
$ws.Cells.Item(2, 36).Value = 15442854

# Row 3
$ws.Cells.Item(3, 4).Value = 2500
$ws.Cells.Item(3, 5).Value = -392
$ws.Cells.Item(3, 6).Value = -392
$ws.Cells.Item(3, 7).Value = -1225
$ws.Cells.Item(3, 8).Value = -1219
$ws.Cells.Item(3, 9).Value = -1069
$ws.Cells.Item(3, 10).Value = -150
$ws.Cells.Item(3, 11).Value = 6305
$ws.Cells.Item(3, 12).Value = 4761
$ws.Cells.Item(3, 13).Value = 1544
$ws.Cells.Item(3, 14).Value = 1248
$ws.Cells.Item(3, 15).Value = 296
$ws.Cells.Item(3, 16).Value = 666
$ws.Cells.Item(3, 17).Value = -21
$ws.Cells.Item(3, 18).Value = 394
$ws.Cells.Item(3, 19).Value = -407
$ws.Cells.Item(3, 20).Value = 31
$ws.Cells.Item(3, 21).Value = -51
$ws.Cells.Item(3, 22).Value = 3074
$ws.Cells.Item(3, 23).Value = -15.7
$ws.Cells.Item(3, 24).Value = -48.75
$ws.Cells.Item(3, 25).Value = -60.78
$ws.Cells.Item(3, 26).Value = -16.93
$ws.Cells.Item(3, 27).Value = 308.3
$ws.Cells.Item(3, 28).Value = -161.72
$ws.Cells.Item(3, 29).Value = -6837
$ws.Cells.Item(3, 30).Value = -0.66
$ws.Cells.Item(3, 31).Value = 7960
$ws.Cells.Item(3, 32).Value = 0.57
$ws.Cells.Item(3, 33).Value = 0
$ws.Cells.Item(3, 34).Value = 0
$ws.Cells.Item(3, 35).Value = 0
$ws.Cells.Item(3, 36).Value = 15709565

# Row 4
$ws.Cells.Item(4, 4).Value = 2865
$ws.Cells.Item(4, 5).Value = -140
$ws.Cells.Item(4, 6).Value = -140
$ws.Cells.Item(4, 7).Value = -725
$ws.Cells.Item(4, 8).Value = -725
$ws.Cells.Item(4, 9).Value = -734
$ws.Cells.Item(4, 10).Value = 8
$ws.Cells.Item(4, 11).Value = 4781
$ws.Cells.Item(4, 12).Value = 3726
$ws.Cells.Item(4, 13).Value = 1055
$ws.Cells.Item(4, 14).Value = 741
$ws.Cells.Item(4, 15).Value = 314
$ws.Cells.Item(4, 16).Value = 666
$ws.Cells.Item(4, 17).Value = -168
$ws.Cells.Item(4, 18).Value = 1036
$ws.Cells.Item(4, 19).Value = -843
$ws.Cells.Item(4, 20).Value = 100
$ws.Cells.Item(4, 21).Value = -269
$ws.Cells.Item(4, 22).Value = 2322
$ws.Cells.Item(4, 23).Value = -4.88
$ws.Cells.Item(4, 24).Value = -25.31
$ws.Cells.Item(4, 25).Value = -73.73999999999999
$ws.Cells.Item(4, 26).Value = -13.08
$ws.Cells.Item(4, 27).Value = 353.29
$ws.Cells.Item(4, 28).Value = -110.81
$ws.Cells.Item(4, 29).Value = -4669
$ws.Cells.Item(4, 30).Value = -0.77
$ws.Cells.Item(4, 31).Value = 4726
$ws.Cells.Item(4, 32).Value = 0.76
$ws.Cells.Item(4, 33).Value = 0
$ws.Cells.Item(4, 34).Value = 0
$ws.Cells.Item(4, 35).Value = 0
$ws.Cells.Item(4, 36).Value = 15709565

# Row 5
$ws.Cells.Item(5, 4).Value = 4316
$ws.Cells.Item(5, 5).Value = 207
$ws.Cells.Item(5, 6).Value = 207
$ws.Cells.Item(5, 7).Value = 244
$ws.Cells.Item(5, 8).Value = 273
$ws.Cells.Item(5, 9).Value = 214
$ws.Cells.Item(5, 10).Value = 59
$ws.Cells.Item(5, 11).Value = 4548
$ws.Cells.Item(5, 12).Value = 2965
$ws.Cells.Item(5, 13).Value = 1584
$ws.Cells.Item(5, 14).Value = 1057
$ws.Cells.Item(5, 15).Value = 527
$ws.Cells.Item(5, 16).Value = 748
$ws.Cells.Item(5, 17).Value = -349
$ws.Cells.Item(5, 18).Value = 501
$ws.Cells.Item(5, 19).Value = -180
$ws.Cells.Item(5, 20).Value = 139
$ws.Cells.Item(5, 21).Value = -488
$ws.Cells.Item(5, 22).Value = 1927
$ws.Cells.Item(5, 23).Value = 4.79
$ws.Cells.Item(5, 24).Value = 6.33
$ws.Cells.Item(5, 25).Value = 23.85
$ws.Cells.Item(5, 26).Value = 5.86
$ws.Cells.Item(5, 27).Value = 187.23
$ws.Cells.Item(5, 28).Value = -67.27
$ws.Cells.Item(5, 29).Value = 1301
$ws.Cells.Item(5, 30).Value = 17.97
$ws.Cells.Item(5, 31).Value = 6004
$ws.Cells.Item(5, 32).Value = 3.9
$ws.Cells.Item(5, 33).Value = 0
$ws.Cells.Item(5, 34).Value = 0
$ws.Cells.Item(5, 35).Value = 0
$ws.Cells.Item(5, 36).Value = 17630764

# Row 6
$ws.Cells.Item(6, 4).Value = 6781
$ws.Cells.Item(6, 5).Value = 62
$ws.Cells.Item(6, 6).Value = 62
$ws.Cells.Item(6, 7).Value = -87
$ws.Cells.Item(6, 8).Value = -88
$ws.Cells.Item(6, 9).Value = -159
$ws.Cells.Item(6, 11).Value = 4930
$ws.Cells.Item(6, 12).Value = 3305
$ws.Cells.Item(6, 13).Value = 1625
$ws.Cells.Item(6, 14).Value = 988
$ws.Cells.Item(6, 16).Value = 792
$ws.Cells.Item(6, 17).Value = 7
$ws.Cells.Item(6, 18).Value = -353
$ws.Cells.Item(6, 19).Value = 404
$ws.Cells.Item(6, 20).Value = 351
$ws.Cells.Item(6, 21).Value = -344
$ws.Cells.Item(6, 22).Value = 2113
$ws.Cells.Item(6, 23).Value = 0.92
$ws.Cells.Item(6, 24).Value = -1.3
$ws.Cells.Item(6, 25).Value = -15.55
$ws.Cells.Item(6, 26).Value = -1.87
$ws.Cells.Item(6, 27).Value = 203.43
$ws.Cells.Item(6, 28).Value = -77.48999999999999
$ws.Cells.Item(6, 29).Value = -869
$ws.Cells.Item(6, 30).Value = -12.4
$ws.Cells.Item(6, 31).Value = 5298
$ws.Cells.Item(6, 32).Value = 2.03
$ws.Cells.Item(6, 33).Value = 0
$ws.Cells.Item(6, 34).Value = 0
$ws.Cells.Item(6, 35).Value = 0
$ws.Cells.Item(6, 36).Value = 18670396

# Remove forecast rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E)) data columns D:AJ,
# leaving only the row number, period label, and period columns (A,B,C).
$ws.Range("D7:AJ9").ClearContents()
